$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Excel constants used below (kept numeric for maximum COM compatibility)
$xlPasteFormats = -4122

# --- Row 14 ("mini T" label row) becomes the separator row that closes off
# the block above it, the same way rows 1, 3, 11 and 23 already are in this
# sheet (ht="15"). Bump its height accordingly. ---
$ws.Rows.Item(14).RowHeight = 15

# --- Row 15 now opens a brand-new bordered block (rows 15-23), so its cells
# need the "top of block" formatting that row 4 already uses for the very
# first block (rows 4-11), instead of the "interior" formatting shared by
# rows 16-23. Copy *only* the formatting from row 4 onto row 15 so the
# values already in row 15 are left untouched. ---
$ws.Range("A4:I4").Copy() | Out-Null
$ws.Range("A15:I15").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = 0

# --- View state: the window was scrolled down and zoomed out a bit, and a
# different cell ended up selected. ---
$ws.Range("E16").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 11
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Zoom = 113
$ws.Range("E16").Select() | Out-Null
